$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.594.22'
$ws.Range('E2').Value = '  +0.00%  '
$ws.Range('D3').Value = '1.895.71'
$ws.Range('E3').Value = '  -0.05%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '246.91'
$ws.Range('E5').Value = '  -0.50%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.693'
$ws.Range('E6').Value = '  +0.03%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '43.20'
$ws.Range('E8').Value = '  -1.62%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '56.69'
$ws.Range('E9').Value = '  +9.11%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.357'
$ws.Range('E10').Value = '  +1.48%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0755'
$ws.Range('E11').Value = '  +1.81%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0984'
$ws.Range('E12').Value = '  +1.37%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '14.59'
$ws.Range('E13').Value = '  +11.44%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.794'
$ws.Range('E14').Value = '  +8.76%  '
$ws.Range('D15').Value = '2.174.26'
$ws.Range('E15').Value = '  +0.04%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '5.04'
$ws.Range('E16').Value = '  +1.72%  '
$ws.Range('D17').Value = '1.904.73'
$ws.Range('E17').Value = '  +0.47%  '
$ws.Range('D18').Value = '35.592.84'
$ws.Range('E18').Value = '  -0.01%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '73.57'
$ws.Range('E19').Value = '  -0.37%  '
$ws.Range('D20').Value = '0.0₃0832'
$ws.Range('E20').Value = '  +0.71%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '247.25'
$ws.Range('E21').Value = '  -0.01%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '13.06'
$ws.Range('E22').Value = '  +1.30%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.19'
$ws.Range('E23').Value = '  +4.50%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.68'
$ws.Range('E24').Value = '  +4.65%  '
$ws.Range('E25').Value = '  +0.01%  '
$ws.Range('E26').Value = '  -1.99%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '166.35'
$ws.Range('E27').Value = '  +0.29%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '8.70'
$ws.Range('E28').Value = '  +2.12%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '18.40'
$ws.Range('E29').Value = '  +0.01%  '
$ws.Range('E30').Value = '  +0.25%  '
$ws.Range('E31').Value = '  +4.06%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.0609'
$ws.Range('E32').Value = '  +4.63%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.27'
$ws.Range('E33').Value = '  +1.12%  '
$ws.Range('B34').Value = 'WEMIXToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.85'
$ws.Range('E34').Value = '  +19.19%  '
$ws.Range('B35').Value = 'BinanceUSD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.48'
$ws.Range('E36').Value = '  -15.82%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.856'
$ws.Range('E37').Value = '  +0.58%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0744'
$ws.Range('E38').Value = '  +8.70%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.95'
$ws.Range('E39').Value = '  -3.40%  '
$ws.Range('E40').Value = '  +7.70%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '99.34'
$ws.Range('E41').Value = '  +1.28%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '16.96'
$ws.Range('E42').Value = '  -1.39%  '
$ws.Range('B43').Value = 'Gas'
$ws.Range('C43').Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '14.48'
$ws.Range('E43').Value = '  +19.79%  '
$ws.Range('B44').Value = 'ARBITRUM'
$ws.Range('C44').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.09'
$ws.Range('E44').Value = '  -0.53%  '
$ws.Range('D45').Value = '1.314.93'
$ws.Range('E45').Value = '  +1.18%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.35'
$ws.Range('E46').Value = '  -1.03%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0813'
$ws.Range('E47').Value = '  +0.37%  '
$ws.Range('E48').Value = '  +0.04%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.74'
$ws.Range('E49').Value = '  -0.39%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '6.38'
$ws.Range('E50').Value = '  +0.04%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '42.63'
$ws.Range('E51').Value = '  -1.96%  '
